$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B. This shifts the existing
# B:K columns (PercActivations..totalStd) to C:L, and leaves column A
# (the segment names) untouched for now - we will replace it below.
$ws.Columns.Item(2).Insert()

# New header for the inserted column B
$ws.Cells.Item(1, 2).Value = "segments"

# Determine last used row (20 in the source data: header + 19 segments)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $segmentName = $ws.Cells.Item($r, 1).Value
    # Move the segment name text into the new column B
    $ws.Cells.Item($r, 2).Value = $segmentName
    # Replace column A with a numeric zero-based index
    $ws.Cells.Item($r, 1).Value = $r - 2
}
